$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 14 data (continuing the log table: "Bài19" / "1. Validate dữ liệu" / commit link) ---
$ws.Range("A14").Value = "Bài19"
$ws.Range("B14").Value = "1. Validate dữ liệu"
$ws.Range("C14").Value = "https://github.com/nguyentienminh07102004/product-management/commit/c718b5b12a9cfe129b0fd827888980b0087cc7d7"

# Match B14's look-and-feel to the rest of the "Bài 18" block (same fill as B12/B13)
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Turn C14 into a hyperlink pointing at the commit, same as the rest of column C
$ws.Hyperlinks.Add($ws.Range("C14"), "https://github.com/nguyentienminh07102004/product-management/commit/c718b5b12a9cfe129b0fd827888980b0087cc7d7") | Out-Null

# --- Selection bookkeeping to mirror what Excel leaves behind after this edit ---
$ws.Range("C15").Select()
